$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.690.60"
$ws.Range("E2").Value = "  +2.05%  "
$ws.Range("D3").Value = "1.869.93"
$ws.Range("E3").Value = "  +1.86%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.008"
$ws.Range("E4").Value = "  +0.53%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "324.66"
$ws.Range("E5").Value = "  -0.21%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.007"
$ws.Range("E6").Value = "  +0.44%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4563"
$ws.Range("E7").Value = "  -1.30%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3837"
$ws.Range("E8").Value = "  -0.63%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07808"
$ws.Range("E9").Value = "  -0.45%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9867"
$ws.Range("E10").Value = "  +2.70%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "21.53"
$ws.Range("E11").Value = "  -2.03%  "
$ws.Range("D12").Value = "1.875.54"
$ws.Range("E12").Value = "  +0.49%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.942"
$ws.Range("E13").Value = "  +0.79%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.645"
$ws.Range("E14").Value = "  -0.46%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.06955"
$ws.Range("E15").Value = "  +1.46%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "87.98"
$ws.Range("E16").Value = "  -0.37%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.008"
$ws.Range("E17").Value = "  +0.54%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000009936"
$ws.Range("E18").Value = "  +0.05%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "16.69"
$ws.Range("E19").Value = "  +0.08%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.006"
$ws.Range("E20").Value = "  +0.40%  "
$ws.Range("D21").Value = "28.710.27"
$ws.Range("E21").Value = "  +2.07%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.248"
$ws.Range("E22").Value = "  -0.84%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.94"
$ws.Range("E23").Value = "  -0.50%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.115"
$ws.Range("E24").Value = "  +0.87%  "
$ws.Range("D25").Value = "2.122.12"
$ws.Range("E25").Value = "  +2.17%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "152.87"
$ws.Range("E26").Value = "  -0.98%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.16"
$ws.Range("E27").Value = "  +0.36%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.713"
$ws.Range("E28").Value = "  +0.26%  "
$ws.Range("B29").Value = "LidoDAOToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.934"
$ws.Range("E29").Value = "  -1.94%  "
$ws.Range("B30").Value = "BitcoinCash"
$ws.Range("C30").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "118.23"
$ws.Range("E30").Value = "  -0.62%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09277"
$ws.Range("E31").Value = "  +0.21%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.9066"
$ws.Range("E32").Value = "  -3.08%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.274"
$ws.Range("E33").Value = "  +0.17%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.317"
$ws.Range("E34").Value = "  -0.19%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.305"
$ws.Range("E35").Value = "  -0.71%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.05701"
$ws.Range("E36").Value = "  -1.97%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.140"
$ws.Range("E37").Value = "  +0.18%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02065"
$ws.Range("E38").Value = "  -2.40%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "7.646"
$ws.Range("E39").Value = "  -0.82%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5599"
$ws.Range("E40").Value = "  +0.17%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1771"
$ws.Range("E41").Value = "  +1.00%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "9.664"
$ws.Range("E42").Value = "  -2.64%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.07097"
$ws.Range("E43").Value = "  -3.47%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "11.56"
$ws.Range("E44").Value = "  -0.80%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.5240"
$ws.Range("E45").Value = "  -0.42%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.121"
$ws.Range("E46").Value = "  -0.67%  "
$ws.Range("B47").Value = "NEARProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.817"
$ws.Range("E47").Value = "  -1.05%  "
$ws.Range("B48").Value = "Quant"
$ws.Range("C48").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "112.74"
$ws.Range("E48").Value = "  -0.71%  "
$ws.Range("B49").Value = "WEMIXToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.093"
$ws.Range("E49").Value = "  -4.26%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.427"
$ws.Range("E50").Value = "  +4.65%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.006"
$ws.Range("E51").Value = "  +0.45%  "

$ws.Range("D4").Style = "Normal"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").Style = "Normal"
